$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the email address stored next to the "Email" label (B7)
$ws.Range("B7").Value = "andra.andruta60@gmail.com"

# Give column B (the Value column) an explicit width, matching the
# resaved layout - closest achievable width under this engine's
# pixel-quantized column-width model is ColumnWidth ~= 24.14
$ws.Columns.Item(2).ColumnWidth = 24.14
